$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25 (shifts existing rows 25-31 down to 26-32)
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new price-report entry
$ws.Cells.Item(25, 1).Value = 5
$ws.Cells.Item(25, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(25, 3).Value = "Maule"
$ws.Cells.Item(25, 4).Value = 44726
$ws.Cells.Item(25, 5).Value = 7
$ws.Cells.Item(25, 6).Value = 100112043
$ws.Cells.Item(25, 7).Value = "Pepino dulce"
$ws.Cells.Item(25, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 300
$ws.Cells.Item(25, 11).Value = 14000
$ws.Cells.Item(25, 12).Value = 14000
$ws.Cells.Item(25, 13).Value = 14000
$ws.Cells.Item(25, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(25, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(25, 16).Value = 778
$ws.Cells.Item(25, 17).Value = 18
$ws.Cells.Item(25, 18).Value = "Hortaliza"
